$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 44
$ws.Range("A44").Value = '2024-10-28 10:12:33'
$ws.Range("B44").Value = 'ORM-0516002'
$ws.Range("C44").Value = 'Z195631098'
$ws.Range("D44").Value = 'C18H19CL2N5OS'
$ws.Range("E44").Value = 'Success'
$ws.Range("F44").Value = '/home/robekott/ERAT/examples/compound_test.sdf'
$ws.Range("G44").Value = 201
$ws.Range("H44").Value = '{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:671f47701a12207a92bcb000"},"data":{"type":"material","id":"asset:671f47701a12207a92bcb000","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:671f47701a12207a92bcb000"},"attributes":{"library":"Compounds","assetTypeId":"5d6e0287ee35880008c18db5","assetId":"671f47701a12207a92bcb000","id":"asset:671f47701a12207a92bcb000","eid":"asset:671f47701a12207a92bcb000","name":"ORM-0516002","synonyms":["CN(CC(=O)NC1C(CL)CCCC1CL)CC1NC(N)C2C(C)C(C)SC2N1","C18H19CL2N5OS"],"description":"","createdAt":"2024-10-28T08:12:33.005Z","editedAt":"2024-10-28T08:12:33.005Z","type":"asset","digest":"16083839","fields":{"Chemical Name":{"value":"2-[({4-amino-5,6-dimethylthieno[2,3-d]pyrimidin-2-yl}methyl)(methyl)amino]-N-(2,6-dichlorophenyl)acetamide"},"Description":{"value":""},"Exact Mass":{"value":"423.06874"},"Material Library Type":{"value":"Compounds"},"Molecular Formula":{"value":"C<sub>18</sub>H<sub>19</sub>Cl<sub>2</sub>N<sub>5</sub>OS"},"Molecular Weight":{"value":"424.34 g/mol"},"Name":{"value":"ORM-0516002"},"Stereochemistry":{"value":"No stereochemistry"}},"flags":{"canTrash":true}},"relationships":{"batches":{"data":[{"type":"material","id":"batch:671f47711a12207a92bcb001","meta":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/batch:671f47711a12207a92bcb001"}}}]},"ancestors":{"data":[{"type":"material","id":"assetType:5d6e0287ee35880008c18db5","meta":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/assetType:5d6e0287ee35880008c18db5"}}}]},"chemicalDrawing":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:671f47701a12207a92bcb000/drawing?format=cdxml"},"data":{"type":"materialDrawing","id":"asset:671f47701a12207a92bcb000"}},"createdBy":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}},"editedBy":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}},"owner":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}}}},"included":[{"type":"material","id":"batch:671f47711a12207a92bcb001","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/batch:671f47711a12207a92bcb001"},"attributes":{"type":"batch","eid":"batch:671f47711a12207a92bcb001","name":"ORM-0516002-001","digest":"73997506","fields":{"Batch Chemical Name":{"value":"2-[({4-amino-5,6-dimethylthieno[2,3-d]pyrimidin-2-yl}methyl)(methyl)amino]-N-(2,6-dichlorophenyl)acetamide"},"Batch Molecular Formula":{"value":"C<sub>18</sub>H<sub>19</sub>Cl<sub>2</sub>N<sub>5</sub>OS"},"Batch Purpose":{"value":"Dummy compound"},"Batch Type":{"value":"Discovery"},"Chemist":{"value":"TestUser MCChemist"},"Description":{"value":""},"Formula Mass":{"value":"424.34 g/mol"},"Name":{"value":"ORM-0516002-001"},"Project":{"value":"Unspecified"},"Source":{"value":"Internal"},"Submission Date":{"value":"2024-10-28T08:12:33.547Z"},"Submitter":{"value":"{userId=114, userName=timo.kangasperko@orion.fi, flags={isSystemStandardUser=true}, alias=timoka, email=timo.kangasperko@orion.fi, firstName=Timo, lastName=Kangasperko, picture={}, isEnabled=true}"},"Synthesis Date":{"value":"2011-10-10T14:48Z"}}}},{"type":"materialDrawing","id":"asset:671f47701a12207a92bcb000","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:671f47701a12207a92bcb000/drawing?format=cdxml"},"attributes":{"id":"asset:671f47701a12207a92bcb000","type":"CHEMICAL_DRAWING"}},{"type":"material","id":"assetType:5d6e0287ee35880008c18db5","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/assetType:5d6e0287ee35880008c18db5"},"attributes":{"type":"assetType","eid":"assetType:5d6e0287ee35880008c18db5","name":"Compounds","digest":"13184553","fields":{}}},{"type":"user","id":"114","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"attributes":{"userId":"114","userName":"timo.kangasperko@orion.fi","flags":{"isSystemStandardUser":true},"alias":"timoka","email":"timo.kangasperko@orion.fi","firstName":"Timo","lastName":"Kangasperko","isEnabled":true},"relationships":{"systemGroups":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114/systemGroups"}}}}]}'

# Row 45
$ws.Range("A45").Value = '2024-10-28 10:12:37'
$ws.Range("B45").Value = 'ORM-0516003'
$ws.Range("C45").Value = 'Z2754556176'
$ws.Range("D45").Value = 'C17H28N4O2'
$ws.Range("E45").Value = 'Success'
$ws.Range("F45").Value = '/home/robekott/ERAT/examples/compound_test.sdf'
$ws.Range("G45").Value = 201
$ws.Range("H45").Value = '{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:671f4774bc105f6067d3957e"},"data":{"type":"material","id":"asset:671f4774bc105f6067d3957e","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:671f4774bc105f6067d3957e"},"attributes":{"library":"Compounds","assetTypeId":"5d6e0287ee35880008c18db5","assetId":"671f4774bc105f6067d3957e","id":"asset:671f4774bc105f6067d3957e","eid":"asset:671f4774bc105f6067d3957e","name":"ORM-0516003","synonyms":["COCCN1CCNC1C1N(CC1(C)C)C(=O)CC1(CN)CC1","C17H28N4O2"],"description":"","createdAt":"2024-10-28T08:12:36.905Z","editedAt":"2024-10-28T08:12:36.905Z","type":"asset","digest":"38794071","fields":{"Chemical Name":{"value":"2-[1-(aminomethyl)cyclopropyl]-1-{2-[1-(2-methoxyethyl)-1H-imidazol-2-yl]-3,3-dimethylazetidin-1-yl}ethan-1-one hydrobromide"},"Description":{"value":""},"Exact Mass":{"value":"320.22123"},"Material Library Type":{"value":"Compounds"},"Molecular Formula":{"value":"C<sub>17</sub>H<sub>28</sub>N<sub>4</sub>O<sub>2</sub>"},"Molecular Weight":{"value":"320.44 g/mol"},"Name":{"value":"ORM-0516003"},"Stereochemistry":{"value":"No stereochemistry"}},"flags":{"canTrash":true}},"relationships":{"batches":{"data":[{"type":"material","id":"batch:671f4775bc105f6067d3957f","meta":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/batch:671f4775bc105f6067d3957f"}}}]},"ancestors":{"data":[{"type":"material","id":"assetType:5d6e0287ee35880008c18db5","meta":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/assetType:5d6e0287ee35880008c18db5"}}}]},"chemicalDrawing":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:671f4774bc105f6067d3957e/drawing?format=cdxml"},"data":{"type":"materialDrawing","id":"asset:671f4774bc105f6067d3957e"}},"createdBy":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}},"editedBy":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}},"owner":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"data":{"type":"user","id":"114"}}}},"included":[{"type":"material","id":"batch:671f4775bc105f6067d3957f","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/batch:671f4775bc105f6067d3957f"},"attributes":{"type":"batch","eid":"batch:671f4775bc105f6067d3957f","name":"ORM-0516003-001","digest":"32351343","fields":{"Batch Chemical Name":{"value":"2-[1-(aminomethyl)cyclopropyl]-1-{2-[1-(2-methoxyethyl)-1H-imidazol-2-yl]-3,3-dimethylazetidin-1-yl}ethan-1-one hydrobromide, hydrogen bromide"},"Batch Molecular Formula":{"value":"C<sub>17</sub>H<sub>28</sub>N<sub>4</sub>O<sub>2</sub>&middot;BRH"},"Batch Purpose":{"value":"Dummy compound"},"Batch Type":{"value":"Discovery"},"Chemist":{"value":"TestUser MCChemist"},"Description":{"value":""},"Formula Mass":{"value":"359.901 g/mol"},"Name":{"value":"ORM-0516003-001"},"Project":{"value":"Unspecified"},"Source":{"value":"Internal"},"Submission Date":{"value":"2024-10-28T08:12:37.534Z"},"Submitter":{"value":"{userId=114, userName=timo.kangasperko@orion.fi, flags={isSystemStandardUser=true}, alias=timoka, email=timo.kangasperko@orion.fi, firstName=Timo, lastName=Kangasperko, picture={}, isEnabled=true}"},"Synthesis Date":{"value":"2011-10-10T14:48Z"}}}},{"type":"material","id":"assetType:5d6e0287ee35880008c18db5","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/assetType:5d6e0287ee35880008c18db5"},"attributes":{"type":"assetType","eid":"assetType:5d6e0287ee35880008c18db5","name":"Compounds","digest":"13184553","fields":{}}},{"type":"user","id":"114","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114"},"attributes":{"userId":"114","userName":"timo.kangasperko@orion.fi","flags":{"isSystemStandardUser":true},"alias":"timoka","email":"timo.kangasperko@orion.fi","firstName":"Timo","lastName":"Kangasperko","isEnabled":true},"relationships":{"systemGroups":{"links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/users/114/systemGroups"}}}},{"type":"materialDrawing","id":"asset:671f4774bc105f6067d3957e","links":{"self":"https://orionsandbox.signalsresearch.revvitycloud.eu/api/rest/v1.0/materials/asset:671f4774bc105f6067d3957e/drawing?format=cdxml"},"attributes":{"id":"asset:671f4774bc105f6067d3957e","type":"CHEMICAL_DRAWING"}}]}'
